# Update "想去人数" (F column) figures across the three sheets that carry
# this data: 展览, 演出, and the combined 全部类型 sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 859
$ws1.Range("F3").Value  = 1434
$ws1.Range("F4").Value  = 1079
$ws1.Range("F5").Value  = 506
$ws1.Range("F6").Value  = 215
$ws1.Range("F7").Value  = 656
$ws1.Range("F8").Value  = 231
$ws1.Range("F10").Value = 70
$ws1.Range("F12").Value = 142
$ws1.Range("F13").Value = 1765
$ws1.Range("F14").Value = 423
$ws1.Range("F15").Value = 38
$ws1.Range("F16").Value = 485
$ws1.Range("F17").Value = 250
$ws1.Range("F18").Value = 407
$ws1.Range("F20").Value = 5
$ws1.Range("F21").Value = 653
$ws1.Range("F22").Value = 44
$ws1.Range("F23").Value = 234
$ws1.Range("F24").Value = 954
$ws1.Range("F26").Value = 1521
$ws1.Range("F27").Value = 265

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 6
$ws2.Range("F8").Value = 279

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 859
$ws4.Range("F4").Value  = 1434
$ws4.Range("F5").Value  = 1079
$ws4.Range("F8").Value  = 506
$ws4.Range("F9").Value  = 215
$ws4.Range("F10").Value = 656
$ws4.Range("F12").Value = 231
$ws4.Range("F14").Value = 70
$ws4.Range("F16").Value = 142
$ws4.Range("F17").Value = 1765
$ws4.Range("F19").Value = 423
$ws4.Range("F20").Value = 38
$ws4.Range("F21").Value = 485
$ws4.Range("F22").Value = 250
$ws4.Range("F23").Value = 407
$ws4.Range("F26").Value = 5
$ws4.Range("F27").Value = 6
$ws4.Range("F28").Value = 279
$ws4.Range("F30").Value = 653
$ws4.Range("F35").Value = 44
$ws4.Range("F36").Value = 234
$ws4.Range("F37").Value = 954
$ws4.Range("F39").Value = 1521
$ws4.Range("F40").Value = 265
